$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.720.83"
$ws.Range("E2").Value = "  +2.41%  "
$ws.Range("D3").Value = "2.221.04"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "'241.00"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").Value = "'0.620"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "'74.76"
$ws.Range("E7").Value = "  +3.71%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  +3.30%  "
$ws.Range("D10").Value = "'41.33"
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("D11").Value = "'0.0931"
$ws.Range("E11").Value = "  -0.66%  "
$ws.Range("D12").Value = "'54.92"
$ws.Range("E12").Value = "  -2.10%  "
$ws.Range("D13").Value = "'6.91"
$ws.Range("E13").Value = "  +1.30%  "
$ws.Range("E14").Value = "  -1.41%  "
$ws.Range("D15").Value = "2.554.17"
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("D16").Value = "'14.67"
$ws.Range("E16").Value = "  +4.51%  "
$ws.Range("D17").Value = "2.221.22"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "'0.802"
$ws.Range("E18").Value = "  -2.36%  "
$ws.Range("D19").Value = "42.573.13"
$ws.Range("E19").Value = "  +2.33%  "
$ws.Range("E20").Value = "  +1.47%  "
$ws.Range("D21").Value = "'70.85"
$ws.Range("E21").Value = "  -0.66%  "
$ws.Range("D22").Value = "'5.94"
$ws.Range("E22").Value = "  -2.32%  "
$ws.Range("D23").Value = "'9.88"
$ws.Range("E23").Value = "  -7.08%  "
$ws.Range("D24").Value = "'229.92"
$ws.Range("E24").Value = "  +1.27%  "
$ws.Range("D25").Value = "'2.14"
$ws.Range("E25").Value = "  +6.78%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "'10.94"
$ws.Range("E27").Value = "  -2.44%  "
$ws.Range("E28").Value = "  -7.11%  "
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "'172.96"
$ws.Range("E30").Value = "  +3.86%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "'36.59"
$ws.Range("E31").Value = "  +20.82%  "
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "'2.09"
$ws.Range("E32").Value = "  -4.89%  "
$ws.Range("D33").Value = "'20.26"
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("D34").Value = "'0.0798"
$ws.Range("E34").Value = "  +1.68%  "
$ws.Range("D35").Value = "'5.33"
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("E36").Value = "  -1.08%  "
$ws.Range("D37").Value = "'0.109"
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("D38").Value = "'4.45"
$ws.Range("E38").Value = "  +5.55%  "
$ws.Range("D39").Value = "'0.0322"
$ws.Range("E39").Value = "  +7.92%  "
$ws.Range("D40").Value = "'12.56"
$ws.Range("E40").Value = "  -2.37%  "
$ws.Range("E41").Value = "  +2.64%  "
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("D43").Value = "'60.49"
$ws.Range("E43").Value = "  -4.25%  "
$ws.Range("D44").Value = "'0.197"
$ws.Range("E44").Value = "  +1.95%  "
$ws.Range("D45").Value = "'8.59"
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("D46").Value = "'0.0991"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").Value = "'99.24"
$ws.Range("E47").Value = "  -1.43%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'1.11"
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("B49").Value = "WOONetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D49").Value = "'0.441"
$ws.Range("E49").Value = "  +21.71%  "
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("E51").Value = "  -1.15%  "